# MasterExecutor_Sanity - add new test case row (TC25_Verify_product_comparison)
# and flip the Runmode on the first row from Yes to No.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at row 14 (everything from the old row 14 down shifts
#    to row+1; dimension/sheetData get updated automatically by the engine).
# ---------------------------------------------------------------------------
$ws.Rows("14:14").Insert()

# ---------------------------------------------------------------------------
# 2. Populate the new row's values.
#    Columns: A=Section/Page  B=Functionality  C=Testcase_number
#             D=Testcase_description  E=Runmode  F=Severity
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "ALL_PAGES"
$ws.Range("B14").Value = "END_TO_END"
$ws.Range("C14").Value = "TC25_Verify_product_comparison"
$ws.Range("D14").Value = "Verifying Product comparison "
$ws.Range("E14").Value = "Yes"
$ws.Range("F14").Value = "High"

# ---------------------------------------------------------------------------
# 3. Copy cell formatting onto the new row from cells that already carry the
#    desired style so we don't mint stray new style entries.
# ---------------------------------------------------------------------------
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B14").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("E14").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D14").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("F14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Flip the Runmode of the very first data row (row 2) from Yes to No.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "No"

# ---------------------------------------------------------------------------
# 5. Grow the hidden _FilterDatabase defined name by one row, matching the
#    inserted row (A1:F31 -> A1:F32).
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "MasterExecutor!_FilterDatabase") {
        $n.RefersTo = "=MasterExecutor!`$A`$1:`$F`$32"
    }
}

# ---------------------------------------------------------------------------
# 6. Re-point the conditional formatting that highlights duplicate Severity
#    values, shifting every row reference down by one to follow the insert.
# ---------------------------------------------------------------------------
$fcs = $ws.Cells.FormatConditions
$fc1 = $fcs.Item(1)
$fc1.ModifyAppliesToRange($ws.Range("F26"))

$fc2 = $fcs.Item(2)
$fc2.ModifyAppliesToRange($ws.Range("F22:F25"))

$fc3 = $ws.Range("F27:F32").FormatConditions.AddUniqueValues()
$fc3.Font.Color = 393372
$fc3.Interior.Color = 13551615

# ---------------------------------------------------------------------------
# 7. Reset the view: select D2 (clears any stale top-left/scroll position
#    and collapses the old multi-cell selection down to a single cell).
# ---------------------------------------------------------------------------
$ws.Range("D2").Select()
